# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# This script updates the DAMSLTag (column I) and DialogAct (column J) values
# for the specific rows identified by the re-annotation pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Row = 5; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 30; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 33; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 34; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 35; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 38; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 41; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 43; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 50; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 74; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 81; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 90; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 114; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 132; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 140; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 151; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 156; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 167; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 177; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 201; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 222; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 227; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 235; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 237; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 253; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 285; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 293; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 301; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 303; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 311; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 327; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 331; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 333; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 338; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 341; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 345; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 356; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 364; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 366; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 369; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 372; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 375; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 394; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 401; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 424; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 499; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 500; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 518; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 522; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 524; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 527; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 530; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 540; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 563; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 564; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 566; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 569; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 571; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 577; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 596; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 598; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 601; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 602; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 604; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 605; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 620; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 628; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 630; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 636; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 637; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 645; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' }
    @{ Row = 653; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 656; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 657; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 679; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 683; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 685; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 686; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 690; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 691; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' }
    @{ Row = 711; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 713; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 714; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 718; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' }
    @{ Row = 727; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
    @{ Row = 730; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
    @{ Row = 734; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 735; DAMSLTag = '%'; DialogAct = 'Uninterpretable' }
    @{ Row = 741; DAMSLTag = 'ba'; DialogAct = 'Appreciation' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.DAMSLTag
    $ws.Cells.Item($u.Row, 10).Value = $u.DialogAct
}

Write-Output "Updated $($updates.Count) rows"
